# Applies the Jenova_Profits market-data refresh across all 8 crafting-class sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 17: One for the Road
$ws.Range("H17").Value = 78945.08
$ws.Range("J17").Value = 78945.08
$ws.Range("L17").Value = 236835.24
$ws.Range("N17").Value = -237171.24

# ALC row 40: Stuck in the Moment
$ws.Range("H40").Value = 6304.923
$ws.Range("I40").Value = 4107.4443
$ws.Range("K40").Value = 4107.4443
$ws.Range("M40").Value = -3932.4443

# ALC row 112: Making Ends Meet
$ws.Range("H112").Value = 3138.932
$ws.Range("J112").Value = 3254.238
$ws.Range("L112").Value = 9762.714
$ws.Range("N112").Value = -11978.714

# ALC row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 3954.889
$ws.Range("I137").Value = 2358.6924
$ws.Range("J137").Value = 6139.1577
$ws.Range("K137").Value = 7076.0772
$ws.Range("L137").Value = 18417.4731
$ws.Range("M137").Value = -4526.0772
$ws.Range("N137").Value = -23517.4731

# ALC row 138: All-night Crafting
$ws.Range("H138").Value = 4230.543
$ws.Range("J138").Value = 5607.0625
$ws.Range("L138").Value = 16821.1875
$ws.Range("N138").Value = -27101.1875

$ws = $wb.Worksheets.Item("ARM")
# ARM row 17: Cook Intentions
$ws.Range("H17").Value = 5599.75
$ws.Range("J17").Value = 5599.75
$ws.Range("L17").Value = 5599.75
$ws.Range("N17").Value = -5945.75

# ARM row 32: Ingot We Trust
$ws.Range("H32").Value = 2877.7415
$ws.Range("I32").Value = 2264.804
$ws.Range("K32").Value = 2264.804
$ws.Range("M32").Value = -1977.804

# ARM row 35: Need for Mead
$ws.Range("H35").Value = 1015
$ws.Range("I35").Value = 1015
$ws.Range("K35").Value = 1015
$ws.Range("M35").Value = -609

# ARM row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 2296.087
$ws.Range("I61").Value = 1540.5
$ws.Range("K61").Value = 1540.5
$ws.Range("M61").Value = -1328.5

# ARM row 74: As the Bolt Flies
$ws.Range("H74").Value = 1411.5476
$ws.Range("J74").Value = 1334.6666
$ws.Range("L74").Value = 1334.6666
$ws.Range("N74").Value = -3082.6666

# ARM row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 1411.5476
$ws.Range("J77").Value = 1334.6666
$ws.Range("L77").Value = 6673.333000000001
$ws.Range("N77").Value = -15409.333

# ARM row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 2997.4385
$ws.Range("I132").Value = 996.8475
$ws.Range("K132").Value = 2990.5425
$ws.Range("M132").Value = -460.5425

# ARM row 136: Metal with Mettle
$ws.Range("H136").Value = 2296.087
$ws.Range("I136").Value = 1540.5
$ws.Range("K136").Value = 4621.5
$ws.Range("M136").Value = -2071.5

$ws = $wb.Worksheets.Item("BSM")
# BSM row 20: Smelt and Dealt
$ws.Range("H20").Value = 3028.7
$ws.Range("I20").Value = 2148.5
$ws.Range("K20").Value = 2148.5
$ws.Range("M20").Value = -1901.5

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31: Wall Not Found
$ws.Range("H31").Value = 559221.8
$ws.Range("I31").Value = 1002149.3
$ws.Range("K31").Value = 1002149.3
$ws.Range("M31").Value = -1001854.3

# CRP row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 559221.8
$ws.Range("I34").Value = 1002149.3
$ws.Range("K34").Value = 1002149.3
$ws.Range("M34").Value = -1001947.3

# CRP row 68: Do You Even String Bow
$ws.Range("H68").Value = 40851.766
$ws.Range("J68").Value = 84542.86
$ws.Range("L68").Value = 84542.86
$ws.Range("N68").Value = -86040.86

# CRP row 71: Win One Bow, Get Three Free (L)
$ws.Range("H71").Value = 40851.766
$ws.Range("J71").Value = 84542.86
$ws.Range("L71").Value = 253628.58
$ws.Range("N71").Value = -261116.58

# CRP row 109: Playing the Market
$ws.Range("H109").Value = 74993.5
$ws.Range("J109").Value = 74993.5
$ws.Range("L109").Value = 74993.5
$ws.Range("N109").Value = -77073.5

# CRP row 121: Safety First
$ws.Range("H121").Value = 99999.5
$ws.Range("J121").Value = 99999.5
$ws.Range("L121").Value = 99999.5
$ws.Range("N121").Value = -102619.5

$ws = $wb.Worksheets.Item("CUL")
# CUL row 4: In Hot Water
$ws.Range("H4").Value = 16224210
$ws.Range("I4").Value = 1412918.8
$ws.Range("K4").Value = 4238756.4
$ws.Range("M4").Value = -4238644.4

# CUL row 98: Sweet Kiss of Death
$ws.Range("H98").Value = 1260.0526
$ws.Range("I98").Value = 1729.6666
$ws.Range("K98").Value = 5188.9998
$ws.Range("M98").Value = -3690.9998

# CUL row 128: A Historical Flavor
$ws.Range("H128").Value = 334996
$ws.Range("I128").Value = 334996
$ws.Range("K128").Value = 1004988
$ws.Range("M128").Value = -1000008

# CUL row 131: The Mountain Steeped
$ws.Range("H131").Value = 4258.905
$ws.Range("J131").Value = 5582.6
$ws.Range("L131").Value = 16747.8
$ws.Range("N131").Value = -26827.8

$ws = $wb.Worksheets.Item("GSM")
# GSM row 18: Gorgeous Gorget
$ws.Range("H18").Value = 1000000000
$ws.Range("J18").Value = 1000000000
$ws.Range("L18").Value = 1000000000
$ws.Range("N18").Value = -1000000586

# GSM row 40: A Little Bird Told Me
$ws.Range("H40").Value = 8775
$ws.Range("I40").Value = 2550
$ws.Range("K40").Value = 2550
$ws.Range("M40").Value = -2399

# GSM row 41: Renascence Man
$ws.Range("H41").Value = 7057
$ws.Range("J41").Value = 7057
$ws.Range("L41").Value = 7057
$ws.Range("N41").Value = -7767

# GSM row 70: Sky Is the Limit
$ws.Range("H70").Value = 90915416
$ws.Range("I70").Value = 6924.625
$ws.Range("J70").Value = 333338080
$ws.Range("K70").Value = 6924.625
$ws.Range("L70").Value = 333338080
$ws.Range("M70").Value = -6654.625
$ws.Range("N70").Value = -333338620

# GSM row 73: Hulls of Broken Dreams (L)
$ws.Range("H73").Value = 90915416
$ws.Range("I73").Value = 6924.625
$ws.Range("J73").Value = 333338080
$ws.Range("K73").Value = 6924.625
$ws.Range("L73").Value = 333338080
$ws.Range("M73").Value = -5988.625
$ws.Range("N73").Value = -333339952

# GSM row 80: Needs More Prayerbell
$ws.Range("H80").Value = 1433301.1
$ws.Range("J80").Value = 1114911.9
$ws.Range("L80").Value = 1114911.9
$ws.Range("N80").Value = -1116907.9

# GSM row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 1433301.1
$ws.Range("J83").Value = 1114911.9
$ws.Range("L83").Value = 5574559.5
$ws.Range("N83").Value = -5584543.5

# GSM row 102: Put the Metal to the Peddle
$ws.Range("H102").Value = 1464.5834
$ws.Range("I102").Value = 1615.7778
$ws.Range("J102").Value = 1011
$ws.Range("K102").Value = 1615.7778
$ws.Range("L102").Value = 1011
$ws.Range("M102").Value = 6.22219999999993
$ws.Range("N102").Value = -4255

# GSM row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 2769.0435
$ws.Range("I122").Value = 2512.7646
$ws.Range("K122").Value = 7538.293799999999
$ws.Range("M122").Value = -5088.293799999999

# GSM row 132: On Board for Lar
$ws.Range("H132").Value = 285037.53
$ws.Range("I132").Value = 305574.7
$ws.Range("K132").Value = 916724.1000000001
$ws.Range("M132").Value = -914194.1000000001

$ws = $wb.Worksheets.Item("LTW")
# LTW row 22: Skin off Their Backs
$ws.Range("H22").Value = 2585.818
$ws.Range("I22").Value = 2029.4
$ws.Range("K22").Value = 2029.4
$ws.Range("M22").Value = -1734.4

# LTW row 27: Fire and Hide
$ws.Range("H27").Value = 2585.818
$ws.Range("I27").Value = 2029.4
$ws.Range("K27").Value = 2029.4
$ws.Range("M27").Value = -1922.4

# LTW row 61: Spelling Me Softly
$ws.Range("H61").Value = 6274.778
$ws.Range("I61").Value = 5412.3335
$ws.Range("K61").Value = 5412.3335
$ws.Range("M61").Value = -5210.3335

# LTW row 68: You Could Say It's a Moving Target
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").ClearContents()
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = 0

# LTW row 71: They Call It Bloody Mary (L)
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").ClearContents()
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = 0

# LTW row 113: Peace in Rest
$ws.Range("H113").Value = 6274.778
$ws.Range("I113").Value = 5412.3335
$ws.Range("K113").Value = 5412.3335
$ws.Range("M113").Value = -3242.3335

# LTW row 122: Hell on Leather
$ws.Range("H122").Value = 1440079.4
$ws.Range("I122").Value = 1265912.6
$ws.Range("J122").Value = 1672301.6
$ws.Range("K122").Value = 3797737.8
$ws.Range("L122").Value = 5016904.800000001
$ws.Range("M122").Value = -3795287.8
$ws.Range("N122").Value = -5021804.800000001

# LTW row 132: Tenets of Tanning
$ws.Range("H132").Value = 4069.375
$ws.Range("I132").Value = 3465.3572
$ws.Range("J132").Value = 4915
$ws.Range("K132").Value = 10396.0716
$ws.Range("L132").Value = 14745
$ws.Range("M132").Value = -7866.071599999999
$ws.Range("N132").Value = -19805

$ws = $wb.Worksheets.Item("WVR")
# WVR row 81: Where the Dragonflies, the Net Catches
$ws.Range("H81").Value = 11304.833
$ws.Range("I81").Value = 1715
$ws.Range("J81").Value = 30484.5
$ws.Range("K81").Value = 3430
$ws.Range("L81").Value = 60969
$ws.Range("M81").Value = -2369
$ws.Range("N81").Value = -63091

# WVR row 84: To Kill a Dragon on Nameday (L)
$ws.Range("H84").Value = 11304.833
$ws.Range("I84").Value = 1715
$ws.Range("J84").Value = 30484.5
$ws.Range("K84").Value = 17150
$ws.Range("L84").Value = 304845
$ws.Range("M84").Value = -11846
$ws.Range("N84").Value = -315453

# WVR row 107: Flax Wax
$ws.Range("H107").Value = 83783.414
$ws.Range("J107").Value = 300
$ws.Range("L107").Value = 900
$ws.Range("N107").Value = -4740

# WVR row 113: A Tender Table
$ws.Range("H113").Value = 2910.2856
$ws.Range("I113").Value = 1846.75
$ws.Range("K113").Value = 5540.25
$ws.Range("M113").Value = -3370.25
